$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '71.937.72'
$ws.Range('E2').Value = '  +3.80%  '
$ws.Range('D3').Value = '3.691.41'
$ws.Range('E3').Value = '  +9.19%  '
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '589.30'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E6').Value = '  +1.15%  '
$ws.Range('D7').Value = '3.684.77'
$ws.Range('E7').Value = '  +9.21%  '
$ws.Range('E8').Value = '  +5.33%  '
$ws.Range('E9').Value = '  +0.09%  '
$ws.Range('E10').Value = '  +3.27%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.614'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +4.67%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '50.00'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +3.47%  '
$ws.Range('E13').Value = '  +1.67%  '
$ws.Range('D14').Value = '4.286.59'
$ws.Range('E14').Value = '  +9.24%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '686.30'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.36%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '8.98'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +4.70%  '
$ws.Range('D17').Value = '3.694.25'
$ws.Range('E17').Value = '  +9.40%  '
$ws.Range('D18').Value = '72.066.82'
$ws.Range('E18').Value = '  +3.80%  '
$ws.Range('E19').Value = '  +2.07%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '18.32'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +3.33%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '11.64'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +3.36%  '
$ws.Range('E22').Value = '  +3.97%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '6.10'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +14.23%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '17.80'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +3.92%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '103.85'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +2.69%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '4.05'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +4.40%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.86'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +6.29%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.25'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +5.63%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '35.63'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +6.67%  '
$ws.Range('E30').Value = '  +5.76%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '7.41'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +6.69%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.19'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +9.77%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '577.93'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +5.39%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '11.35'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +3.01%  '
$ws.Range('E35').Value = '  +4.09%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '59.98'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +3.86%  '
$ws.Range('D37').Value = '3.761.35'
$ws.Range('E37').Value = '  +4.37%  '
$ws.Range('E38').Value = '  +0.01%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.145'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +2.93%  '
$ws.Range('D40').Value = '0.0₃0776'
$ws.Range('E40').Value = '  +5.72%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '35.59'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +1.04%  '
$ws.Range('E42').Value = '  +5.15%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.80'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +3.68%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0462'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +8.90%  '
$ws.Range('E45').Value = '  +4.05%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.38'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.36%  '
$ws.Range('E47').Value = '  +6.91%  '
$ws.Range('E48').Value = '  +4.15%  '
$ws.Range('E49').Value = '  +3.89%  '
$ws.Range('E50').Value = '  -0.14%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '134.18'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +3.80%  '
